$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 763.3396
$ws.Range("J17").Value = 763.3396
$ws.Range("L17").Value = 2290.0188
$ws.Range("N17").Value = -2626.0188

$ws.Range("H46").Value = 8595
$ws.Range("I46").Value = 4940
$ws.Range("K46").Value = 14820
$ws.Range("M46").Value = -14701

$ws.Range("H53").Value = 481.4
$ws.Range("J53").Value = 563.6667
$ws.Range("L53").Value = 563.6667
$ws.Range("N53").Value = -1837.6667

$ws.Range("H60").Value = 8595
$ws.Range("I60").Value = 4940
$ws.Range("K60").Value = 14820
$ws.Range("M60").Value = -14336

$ws.Range("H94").Value = 8554974
$ws.Range("I94").Value = 11115566
$ws.Range("K94").Value = 11115566
$ws.Range("M94").Value = -11115115

$ws.Range("H107").Value = 648.0769
$ws.Range("I107").Value = 579.25
$ws.Range("J107").Value = 758.2
$ws.Range("K107").Value = 579.25
$ws.Range("L107").Value = 758.2
$ws.Range("M107").Value = 1340.75
$ws.Range("N107").Value = -4598.2

$ws.Range("H137").Value = 2763.7273
$ws.Range("J137").Value = 3780.4
$ws.Range("L137").Value = 11341.2
$ws.Range("N137").Value = -16441.2

$ws.Range("H138").Value = 5646.9243
$ws.Range("J138").Value = 7254.15
$ws.Range("L138").Value = 21762.45
$ws.Range("N138").Value = -32042.45

$ws.Range("H140").Value = 59164.9
$ws.Range("J140").Value = 59164.9
$ws.Range("L140").Value = 59164.9
$ws.Range("N140").Value = -69524.89999999999

$ws.Range("H141").Value = 4265.1377
$ws.Range("I141").Value = 4449.7827
$ws.Range("J141").Value = 3557.3333
$ws.Range("K141").Value = 13349.3481
$ws.Range("L141").Value = 10671.9999
$ws.Range("M141").Value = -8169.348099999999
$ws.Range("N141").Value = -21031.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2428.0193
$ws.Range("I32").Value = 2368.7144
$ws.Range("J32").Value = 3396.6667
$ws.Range("K32").Value = 2368.7144
$ws.Range("L32").Value = 3396.6667
$ws.Range("M32").Value = -2081.7144
$ws.Range("N32").Value = -3970.6667

$ws.Range("H63").Value = 3641.8333
$ws.Range("I63").Value = 3109.8572
$ws.Range("J63").Value = 4386.6
$ws.Range("K63").Value = 3109.8572
$ws.Range("L63").Value = 4386.6
$ws.Range("M63").Value = -2423.8572
$ws.Range("N63").Value = -5758.6

$ws.Range("H66").Value = 3641.8333
$ws.Range("I66").Value = 3109.8572
$ws.Range("J66").Value = 4386.6
$ws.Range("K66").Value = 15549.286
$ws.Range("L66").Value = 21933
$ws.Range("M66").Value = -12117.286
$ws.Range("N66").Value = -28797

$ws.Range("H74").Value = 13163221
$ws.Range("I74").Value = 41667450
$ws.Range("K74").Value = 41667450
$ws.Range("M74").Value = -41666576

$ws.Range("H77").Value = 13163221
$ws.Range("I77").Value = 41667450
$ws.Range("K77").Value = 208337250
$ws.Range("M77").Value = -208332882

$ws.Range("H88").Value = 150000
$ws.Range("J88").Value = 150000
$ws.Range("L88").Value = 150000
$ws.Range("N88").Value = -150812

$ws.Range("H91").Value = 150000
$ws.Range("J91").Value = 150000
$ws.Range("L91").Value = 150000
$ws.Range("N91").Value = -152808

$ws.Range("H132").Value = 11181.525
$ws.Range("I132").Value = 11384.386
$ws.Range("K132").Value = 34153.158
$ws.Range("M132").Value = -31623.158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2778398.8
$ws.Range("I64").Value = 5208906
$ws.Range("J64").Value = 676.4286
$ws.Range("K64").Value = 5208906
$ws.Range("L64").Value = 676.4286
$ws.Range("M64").Value = -5208681
$ws.Range("N64").Value = -1126.4286

$ws.Range("H67").Value = 2778398.8
$ws.Range("I67").Value = 5208906
$ws.Range("J67").Value = 676.4286
$ws.Range("K67").Value = 5208906
$ws.Range("L67").Value = 676.4286
$ws.Range("M67").Value = -5208126
$ws.Range("N67").Value = -2236.4286

$ws.Range("H86").Value = 2367.75
$ws.Range("I86").Value = 2150
$ws.Range("K86").Value = 2150
$ws.Range("M86").Value = -1027

$ws.Range("H89").Value = 2367.75
$ws.Range("I89").Value = 2150
$ws.Range("K89").Value = 10750
$ws.Range("M89").Value = -5134

$ws.Range("H94").Value = 527365.4
$ws.Range("I94").Value = 596045.2
$ws.Range("J94").Value = 819.6667
$ws.Range("K94").Value = 596045.2
$ws.Range("L94").Value = 819.6667
$ws.Range("M94").Value = -595594.2
$ws.Range("N94").Value = -1721.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23814234
$ws.Range("I31").Value = 71430640
$ws.Range("J31").Value = 6030.0356
$ws.Range("K31").Value = 71430640
$ws.Range("L31").Value = 6030.0356
$ws.Range("M31").Value = -71430345
$ws.Range("N31").Value = -6620.0356

$ws.Range("H34").Value = 23814234
$ws.Range("I34").Value = 71430640
$ws.Range("J34").Value = 6030.0356
$ws.Range("K34").Value = 71430640
$ws.Range("L34").Value = 6030.0356
$ws.Range("M34").Value = -71430438
$ws.Range("N34").Value = -6434.0356

$ws.Range("H99").Value = 13043.667
$ws.Range("I99").Value = 15405.333
$ws.Range("K99").Value = 15405.333
$ws.Range("M99").Value = -13907.333

$ws.Range("H107").Value = 1072863.5
$ws.Range("J107").Value = 5282.857
$ws.Range("L107").Value = 5282.857
$ws.Range("N107").Value = -9122.857

$ws.Range("H126").Value = 13043.667
$ws.Range("I126").Value = 15405.333
$ws.Range("K126").Value = 46215.999
$ws.Range("M126").Value = -43745.999

$ws.Range("H141").Value = 127346.55
$ws.Range("J141").Value = 127346.55
$ws.Range("L141").Value = 127346.55
$ws.Range("N141").Value = -137706.55

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 179992.94
$ws.Range("I68").Value = 1816.5
$ws.Range("K68").Value = 5449.5
$ws.Range("M68").Value = -4638.5

$ws.Range("H71").Value = 179992.94
$ws.Range("I71").Value = 1816.5
$ws.Range("K71").Value = 16348.5
$ws.Range("M71").Value = -12292.5

$ws.Range("H109").Value = 7163.385
$ws.Range("I109").Value = 2790.5
$ws.Range("J109").Value = 14160
$ws.Range("K109").Value = 8371.5
$ws.Range("L109").Value = 42480
$ws.Range("M109").Value = -7331.5
$ws.Range("N109").Value = -44560

$ws.Range("H137").Value = 43002796
$ws.Range("I137").Value = 44118584
$ws.Range("K137").Value = 132355752
$ws.Range("M137").Value = -132350652

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 836718.9
$ws.Range("I80").Value = 1282205.2
$ws.Range("J80").Value = 9387
$ws.Range("K80").Value = 1282205.2
$ws.Range("L80").Value = 9387
$ws.Range("M80").Value = -1281207.2
$ws.Range("N80").Value = -11383

$ws.Range("H83").Value = 836718.9
$ws.Range("I83").Value = 1282205.2
$ws.Range("J83").Value = 9387
$ws.Range("K83").Value = 6411026
$ws.Range("L83").Value = 46935
$ws.Range("M83").Value = -6406034
$ws.Range("N83").Value = -56919

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1749840.2
$ws.Range("I68").Value = 2274312.2
$ws.Range("J68").Value = 1600
$ws.Range("K68").Value = 2274312.2
$ws.Range("L68").Value = 1600
$ws.Range("M68").Value = -2273563.2
$ws.Range("N68").Value = -3098

$ws.Range("H71").Value = 1749840.2
$ws.Range("I71").Value = 2274312.2
$ws.Range("J71").Value = 1600
$ws.Range("K71").Value = 11371561
$ws.Range("L71").Value = 8000
$ws.Range("M71").Value = -11367817
$ws.Range("N71").Value = -15488

$ws.Range("H122").Value = 62505050
$ws.Range("I122").Value = 111115630
$ws.Range("J122").Value = 5727.857
$ws.Range("K122").Value = 333346890
$ws.Range("L122").Value = 17183.571
$ws.Range("M122").Value = -333344440
$ws.Range("N122").Value = -22083.571

$ws.Range("H136").Value = 3925.14
$ws.Range("I136").Value = 3928.7654
$ws.Range("K136").Value = 11786.2962
$ws.Range("M136").Value = -9236.296200000001

$ws.Range("H140").Value = 93889.57000000001
$ws.Range("J140").Value = 93889.57000000001
$ws.Range("L140").Value = 93889.57000000001
$ws.Range("N140").Value = -104249.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 11016.833
$ws.Range("J74").Value = 12020.2
$ws.Range("L74").Value = 12020.2
$ws.Range("N74").Value = -13892.2

$ws.Range("H77").Value = 11016.833
$ws.Range("J77").Value = 12020.2
$ws.Range("L77").Value = 36060.60000000001
$ws.Range("N77").Value = -45420.60000000001

$ws.Range("H107").Value = 1002.7308
$ws.Range("I107").Value = 753.6667
$ws.Range("J107").Value = 2048.8
$ws.Range("K107").Value = 2261.0001
$ws.Range("L107").Value = 6146.400000000001
$ws.Range("M107").Value = -341.0001000000002
$ws.Range("N107").Value = -9986.400000000001

$ws.Range("H113").Value = 674.8095
$ws.Range("J113").Value = 911.2222
$ws.Range("L113").Value = 2733.6666
$ws.Range("N113").Value = -7073.6666

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H136").Value = 5794.15
$ws.Range("I136").Value = 2507.7046
$ws.Range("J136").Value = 8376.357
$ws.Range("K136").Value = 7523.1138
$ws.Range("L136").Value = 25129.071
$ws.Range("M136").Value = -4973.1138
$ws.Range("N136").Value = -30229.071
